# Add a new worksheet "Sheet2" with a small labeled table (x-layout extraction)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Column headers
$ws2.Range("B1").Value = "a"
$ws2.Range("C1").Value = "b"
$ws2.Range("D1").Value = "c"
$ws2.Range("E1").Value = "d"

# Row labels + values
$labels = @("e", "f", "g", "f", "g", "f")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $labels[$i]
    $ws2.Cells.Item($row, 2).Value = 1
    $ws2.Cells.Item($row, 3).Value = 2
    $ws2.Cells.Item($row, 4).Value = 3
    $ws2.Cells.Item($row, 5).Value = 4
}

# Selections
$ws1.Range("A4").Select()
$ws2.Range("B3").Select()

# Make Sheet2 the active sheet/tab
$ws2.Activate()
